# Adding Administrator User Story
#
# Fills in the previously-empty sub-bullet under
# "As an Administrator, I want to search volunteers by last name." and
# appends two additional sub-bullets, all styled like their sibling
# bullets (Garamond font, list level 2 / ilvl=1, numId=2).

$d = $word.ActiveDocument

function Set-BulletText($paragraphIndex, $text) {
    # The target paragraph currently has a single, textless run (just the
    # paragraph mark). Give it a one-character placeholder first so that
    # we can address that run's exact character range, then swap the
    # whole run (text + rPr) via InsertXML. This replaces the run
    # in-place without touching/regenerating the paragraph's own <w:pPr>,
    # so existing paragraph formatting survives untouched.
    $para = $d.Paragraphs.Item($paragraphIndex)
    $para.Range.Text = "X"

    $para = $d.Paragraphs.Item($paragraphIndex)
    $runStart = $para.Range.Start
    $runRange = $d.Range($runStart, $runStart + 1)

    $escaped = $text -replace '&', '&amp;' -replace '<', '&lt;' -replace '>', '&gt;'

    $xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
           '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
           '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
           '<pkg:xmlData>' +
           '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
           '<w:body><w:p><w:r><w:rPr>' +
           '<w:rFonts w:ascii="Garamond" w:cs="Garamond" w:eastAsia="Garamond" w:hAnsi="Garamond"/>' +
           '<w:rtl w:val="0"/>' +
           '</w:rPr><w:t xml:space="preserve">' + $escaped + '</w:t></w:r></w:p></w:body>' +
           '</w:document></pkg:xmlData></pkg:part></pkg:package>'

    $runRange.InsertXML($xml)
}

$bodyText1 = "When the user is prompted by the AdministratorUI, they can choose option 2, which sends them to AdministrarorUI.searchByLastName(), which calls AdministrarorUI.promptForVolsLastName. The last name is then passed to Administrator.getMatchingVolunteers, which returns the list of users with that last name. If there are none, there is a message for that. If there are some, the AdministrarorUI.displayVolunteers() method is called to actually print all of the Volunteers and their information."
$bodyText2 = "There are no error checks associated with this business rule, since the only thing involved is getting a list of Volunteers. Whether the list is populated or not, that is a valid result. "
$bodyText3 = "The volunteers are printed out in alphabetical order, along with what jobs (if any) they are signed up for."

# Locate the existing empty sub-bullet that follows the Administrator user
# story bullet.
$found = $false
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -match "As an Administrator, I want to search volunteers by last name\.") {
        $targetIndex = $i + 1
        $found = $true
    }
}

if (-not $found) {
    throw "Could not find the Administrator user story paragraph."
}

# 1) Fill in the existing (empty) sub-bullet with the first explanatory
#    paragraph.
Set-BulletText $targetIndex $bodyText1

# 2) Clone that paragraph (preserving its list/indent formatting) to make
#    a second sub-bullet immediately after it, then fill it in.
$p = $d.Paragraphs.Item($targetIndex)
$p.Range.InsertParagraphAfter()
Set-BulletText ($targetIndex + 1) $bodyText2

# 3) Clone again for the third sub-bullet.
$p = $d.Paragraphs.Item($targetIndex + 1)
$p.Range.InsertParagraphAfter()
Set-BulletText ($targetIndex + 2) $bodyText3
